# Updates cryptos.xlsx price/volume figures (and two row swaps) per the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Prefix with an apostrophe so Excel stores the value as text even when
    # it looks numeric (e.g. "559.78", "64.620.21"), then drop back to the
    # Normal style so no stray number-format/quote-prefix style lingers on
    # the cell.
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2: Price, Volume(1h)
Set-TextCell "D2" "64.620.21"
Set-TextCell "E2" "  +0.54%  "

# Row 3: Price, Volume(1h)
Set-TextCell "D3" "3.365.93"
Set-TextCell "E3" "  -1.24%  "

# Row 4: Volume(1h)
Set-TextCell "E4" "  -0.07%  "

# Row 5: Price, Volume(1h)
Set-TextCell "D5" "559.78"
Set-TextCell "E5" "  -0.01%  "

# Row 6: Price, Volume(1h)
Set-TextCell "D6" "176.66"
Set-TextCell "E6" "  +2.35%  "

# Row 7: Price, Volume(1h)
Set-TextCell "D7" "0.621"
Set-TextCell "E7" "  +0.63%  "

# Row 8: Price, Volume(1h)
Set-TextCell "D8" "3.355.87"
Set-TextCell "E8" "  -1.35%  "

# Row 9: Volume(1h)
Set-TextCell "E9" "  -0.02%  "

# Row 10: Price, Volume(1h)
Set-TextCell "D10" "0.165"
Set-TextCell "E10" "  +8.39%  "

# Row 11: Price, Volume(1h)
Set-TextCell "D11" "0.631"
Set-TextCell "E11" "  +2.24%  "

# Row 12: Price, Volume(1h)
Set-TextCell "D12" "55.03"
Set-TextCell "E12" "  -1.13%  "

# Row 13: Price, Volume(1h)
Set-TextCell "D13" "0.0000277"
Set-TextCell "E13" "  +3.36%  "

# Row 14: Price
Set-TextCell "D14" "9.12"

# Row 15: Price, Volume(1h)
Set-TextCell "D15" "3.906.91"
Set-TextCell "E15" "  -1.29%  "

# Row 16: Price, Volume(1h)
Set-TextCell "D16" "18.32"
Set-TextCell "E16" "  +2.24%  "

# Row 17: Coin, Link, Price, Volume(1h)
Set-TextCell "B17" "WrappedEther"
Set-TextCell "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D17" "3.373.47"
Set-TextCell "E17" "  -1.04%  "

# Row 18: Coin, Link, Price, Volume(1h)
Set-TextCell "B18" "TRON"
Set-TextCell "C18" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell "D18" "0.117"
Set-TextCell "E18" "  -1.99%  "

# Row 19: Price, Volume(1h)
Set-TextCell "D19" "11.87"
Set-TextCell "E19" "  +1.37%  "

# Row 20: Price, Volume(1h)
Set-TextCell "D20" "64.544.22"
Set-TextCell "E20" "  +0.43%  "

# Row 21: Price, Volume(1h)
Set-TextCell "D21" "0.989"
Set-TextCell "E21" "  +0.37%  "

# Row 22: Price, Volume(1h)
Set-TextCell "D22" "460.02"
Set-TextCell "E22" "  +13.22%  "

# Row 23: Price, Volume(1h)
Set-TextCell "D23" "4.79"
Set-TextCell "E23" "  +10.95%  "

# Row 24: Price, Volume(1h)
Set-TextCell "D24" "4.10"
Set-TextCell "E24" "  -0.08%  "

# Row 25: Price, Volume(1h)
Set-TextCell "D25" "86.08"
Set-TextCell "E25" "  +4.17%  "

# Row 26: Volume(1h)
Set-TextCell "E26" "  +2.27%  "

# Row 27: Price, Volume(1h)
Set-TextCell "D27" "10.89"
Set-TextCell "E27" "  +1.74%  "

# Row 28: Price, Volume(1h)
Set-TextCell "D28" "2.85"
Set-TextCell "E28" "  +3.82%  "

# Row 29: Price, Volume(1h)
Set-TextCell "D29" "8.80"
Set-TextCell "E29" "  +0.15%  "

# Row 30: Price, Volume(1h)
Set-TextCell "D30" "30.15"
Set-TextCell "E30" "  +2.16%  "

# Row 31: Price, Volume(1h)
Set-TextCell "D31" "6.77"
Set-TextCell "E31" "  +3.34%  "

# Row 32: Price, Volume(1h)
Set-TextCell "D32" "11.49"
Set-TextCell "E32" "  +0.70%  "

# Row 33: Price, Volume(1h)
Set-TextCell "D33" "579.36"
Set-TextCell "E33" "  -1.81%  "

# Row 34: Price, Volume(1h)
Set-TextCell "D34" "0.108"
Set-TextCell "E34" "  +1.66%  "

# Row 35: Price, Volume(1h)
Set-TextCell "D35" "59.42"
Set-TextCell "E35" "  +0.88%  "

# Row 36: Volume(1h)
Set-TextCell "E36" "  -0.10%  "

# Row 37: Price, Volume(1h)
Set-TextCell "D37" "0.141"
Set-TextCell "E37" "  -7.23%  "

# Row 38: Coin, Link, Price, Volume(1h)
Set-TextCell "B38" "InjectiveProtocol"
Set-TextCell "C38" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D38" "36.03"
Set-TextCell "E38" "  +1.34%  "

# Row 39: Coin, Link, Price, Volume(1h)
Set-TextCell "B39" "PEPE"
Set-TextCell "C39" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell "D39" "0.0₃0760"
Set-TextCell "E39" "  +4.12%  "

# Row 40: Price, Volume(1h)
Set-TextCell "D40" "3.46"
Set-TextCell "E40" "  +2.20%  "

# Row 41: Price, Volume(1h)
Set-TextCell "D41" "0.373"
Set-TextCell "E41" "  +0.93%  "

# Row 42: Price, Volume(1h)
Set-TextCell "D42" "3.098.21"
Set-TextCell "E42" "  -2.46%  "

# Row 43: Volume(1h)
Set-TextCell "E43" "  -0.05%  "

# Row 44: Price, Volume(1h)
Set-TextCell "D44" "2.85"
Set-TextCell "E44" "  -0.36%  "

# Row 45: Volume(1h)
Set-TextCell "E45" "  +0.72%  "

# Row 46: Price, Volume(1h)
Set-TextCell "D46" "0.0413"
Set-TextCell "E46" "  +2.05%  "

# Row 47: Price, Volume(1h)
Set-TextCell "D47" "3.20"
Set-TextCell "E47" "  -1.24%  "

# Row 48: Price, Volume(1h)
Set-TextCell "D48" "0.131"
Set-TextCell "E48" "  +2.01%  "

# Row 49: Volume(1h)
Set-TextCell "E49" "  -0.87%  "

# Row 50: Price, Volume(1h)
Set-TextCell "D50" "8.37"
Set-TextCell "E50" "  +1.66%  "

# Row 51: Price, Volume(1h)
Set-TextCell "D51" "135.81"
Set-TextCell "E51" "  +1.00%  "

